$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("D3").Interior.Color = 5000268
$ws.Range("E3").Borders.LineStyle = 1
$ws.Range("E3").Interior.Color = 5000268
Write-Host "done"
